# Updated comments to archive jobs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark rows 2-6 (Comment column R) as "Archive"
$ws.Range("R2").Value = "Archive"
$ws.Range("R3").Value = "Archive"
$ws.Range("R4").Value = "Archive"
$ws.Range("R5").Value = "Archive"
$ws.Range("R6").Value = "Archive"

# Flag row 11 for Dee to review
$ws.Range("R11").Value = "Dee to review"

# Move the active selection to R16 to match where the author left off
$ws.Range("R16").Select()
